$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.830.26"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "1.618.22"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.98%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.32"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3934"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.354"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.33"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08452"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.70"
$ws.Range("E13").Value = "  -4.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.047"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.560"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001281"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "1.613.34"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.66"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06929"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.95"
$ws.Range("E20").Value = "  -5.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.812"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.42"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").Value = "23.837.85"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.459"
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.823"
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.23"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.00"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "140.49"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.304"
$ws.Range("E30").Value = "  -8.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.794"
$ws.Range("E31").Value = "  -5.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.495"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").Value = "1.789.85"
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08121"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9830"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.604"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("E37").Value = "  -4.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2666"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09135"
$ws.Range("E39").Value = "  -4.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.31"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.61"
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.427"
$ws.Range("E42").Value = "  -4.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7504"
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.94"
$ws.Range("E44").Value = "  -2.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6916"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.470"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.071"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08243"
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.10"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("E51").Value = "  -8.34%  "
